$wb = $excel.ActiveWorkbook

function Set-TextLabel {
    param($ws, [string]$cellAddr, [string]$text, [string]$scratchAddr = "A1")

    # Directly assigning a numeric-looking string to Range.Value causes Excel
    # to coerce it into a number, which would store it as t="n" instead of
    # text. To force a genuine text value (t="inlineStr"/shared string) while
    # leaving the destination cell's existing style index untouched, stage
    # the text in a scratch cell (using the quote-prefix text-entry trick),
    # copy *values only* into the destination (PasteSpecial values preserves
    # the destination's own style/format), then fully clear the scratch cell
    # so it leaves no residue behind.
    $scratch = $ws.Range($scratchAddr)
    $scratch.Value = "'" + $text
    $scratch.Copy()
    $dest = $ws.Range($cellAddr)
    $dest.PasteSpecial(-4163) # xlPasteValues
    $scratch.Clear()
}

# Sheet 1: "Potencia Acumulada - SIN (MW)"
# (use A13 as scratch -- it is part of the "Total" row being deleted below,
# so any leftover style churn from the text-entry trick disappears with it)
$ws1 = $wb.Worksheets.Item("Potencia Acumulada - SIN (MW)")
Set-TextLabel $ws1 "E1" "2050" "A13"
$ws1.Rows.Item(13).Delete()

# Sheet 2: "Geracao Periodo Medio (MWMed)"
$ws2 = $wb.Worksheets.Item("Geracao Periodo Medio (MWMed)")
Set-TextLabel $ws2 "E1" "2050" "A13"
$ws2.Rows.Item(13).Delete()

# Sheet 3: "Atendimento a Ponta(MW)"
$ws3 = $wb.Worksheets.Item("Atendimento a Ponta(MW)")
Set-TextLabel $ws3 "E1" "2050" "A13"
$ws3.Rows.Item(13).Delete()

# Sheet 4: "Potencia Incremental - SIN(MW)" (headers are ranges, e.g. 2015-2030, 2031-2040)
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
Set-TextLabel $ws4 "E1" "2041-2050" "A13"
$ws4.Rows.Item(13).Delete()

# Sheet 5: "Emissoes Totais (MtCO2eq)" -- only label fix, no Total row present
$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
Set-TextLabel $ws5 "E1" "2050" "A1"

# Sheet 6: "Custo Total (bilhões de R$)" -- remove Total row
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Rows.Item(4).Delete()
